# Remove the <w:contextualSpacing w:val="0"/> element from every paragraph's
# properties (w:pPr), matching the commit's canonical-OOXML diff which drops
# that element from each affected paragraph while leaving everything else
# (borders, shading, run formatting, etc.) untouched.
$d = $word.ActiveDocument

# ContextualSpacing is not exposed as a Word object-model property in this
# host, and Range.WordOpenXML is read-only, so pull the live package XML
# once, compute each paragraph's corrected markup, then push it back in
# place with Range.InsertXML (which replaces only that paragraph's range).
$xml = $d.Content.WordOpenXML

$bodyMatch = [regex]::Match($xml, "<w:body>(.*)</w:body>")
$body = $bodyMatch.Groups[1].Value

$paraMatches = [regex]::Matches($body, "<w:p(?:\s[^>]*)?>.*?</w:p>")

$target = '<w:contextualSpacing w:val="0"/>'
$paraCount = $d.Paragraphs.Count
$changed = 0

for ($i = 0; $i -lt $paraMatches.Count -and $i -lt $paraCount; $i++) {
    $frag = $paraMatches[$i].Value
    if ($frag.Contains($target)) {
        $newFrag = $frag.Replace($target, "")
        $d.Paragraphs($i + 1).Range.InsertXML($newFrag)
        $changed = $changed + 1
    }
}

Write-Output ("Paragraphs updated: " + $changed)
